$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: territorial scope labels (capitalized and relabeled)
$ws.Range("A1").Value = "Municipio"
$ws.Range("B1").Value = "Entidad singular"
$ws.Range("C1").Value = "Núcleo"
$ws.Range("D1").Value = "Comarca"
$ws.Range("E1").Value = "Provincia"

# Row 2: matching iaest-measure identifiers for each territorial scope
$ws.Range("A2").Value = "iaest-measure:municipio"
$ws.Range("B2").Value = "iaest-measure:entidad-singular"
$ws.Range("C2").Value = "iaest-measure:nucleo"
$ws.Range("D2").Value = "iaest-measure:comarca"
$ws.Range("E2").Value = "iaest-measure:provincia"

# Row 4: data type changed from xsd:int to xsd:string
$ws.Range("A4:E4").Value = "xsd:string"
